$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update contact / name cells to include affiliation ---
$ws.Range("G5").Value  = "Jakob Gärtner, DB (LEA)"

$ws.Range("E14").Value = "SQS Team"
$ws.Range("F14").Value = "Bend Hekele, DB"
$ws.Range("G14").Value = "Fausto Cochetti, Alstom"
$ws.Range("J14").Value = "Jan Welvaarts und Vincent Nuhaan, NS Team"
$ws.Range("K14").Value = "Peyman Farhangi, DB"
$ws.Range("L14").Value = "Jakob Gärtner, DB (LEA)"

$ws.Range("F15").Value = "Jakob Gärtner, DB (LEA)"
$ws.Range("G15").Value = "Christian Giraud, Alstom"
$ws.Range("H15").Value = "Alexander Stante, Fraunhofer"
$ws.Range("J15").Value = "Uwe Steinke, Siemens"

$ws.Range("F16").Value = "David Mentre, Mitsubishi"
$ws.Range("G16").Value = "Veronique Gontier, All4Tech"
$ws.Range("J16").ClearContents()

$ws.Range("F17").Value = "Uwe Steinke, Siemens"
$ws.Range("G17").Value = "Benjamin Beichler, UOR"

# --- Column width tweaks ---
# (ColumnWidth is quantized by the engine to 1/6-character pixel steps, so these
#  values are chosen to land on the raw widths closest to the authored
#  27.28515625 / 43.42578125 character-width units.)
$ws.Columns.Item(8).ColumnWidth = 26.5
$ws.Columns.Item(10).ColumnWidth = 42.65

# --- Selection / view state ---
[void]$ws.Range("E18").Select()
